# Add a new "2020" column (column Q) to the indicators table on the
# active sheet, mirroring the formatting of the existing "2019" column (P),
# and update the sheet's current selection, as in the source commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Copy the formatting of column P (the "2019" column, rows 4-14, which
#    includes the header style, the data-row styles and the bottom-border
#    total-row style) onto the new column Q so the new column matches the
#    look of the rest of the table.
$ws.Range("P4:P14").Copy()
$ws.Range("Q4:Q14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 2) Fill in the 2020 values.
$ws.Range("Q4").Value = 2020

$ws.Range("Q5").Value = 0.02
$ws.Range("Q6").Value = 0
$ws.Range("Q7").Value = 0
$ws.Range("Q8").Value = 0
$ws.Range("Q9").Value = 0.54
$ws.Range("Q10").Value = 0
$ws.Range("Q11").Value = 0
$ws.Range("Q12").Value = 0
$ws.Range("Q13").Value = 0
$ws.Range("Q14").Value = 0

# 3) Update the sheet's stored selection/active cell (as saved in the
#    workbook by Excel when it was last edited).
[void]$ws.Range("N19").Select()
